$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The workbook originally held yearly (2008年-2020年) rows of data in rows 2-14.
# The update drops the two oldest years (2008年, 2009年) - shifting all remaining
# rows up by two - and appends a new row for 2021年 at the end.

# Remove the 2008年 row (current row 2); everything below shifts up.
$ws.Rows.Item(2).Delete()
# Remove the 2009年 row (now also row 2, since the previous delete shifted rows up).
$ws.Rows.Item(2).Delete()

# After the two deletions the sheet has 12 data rows (2010年..2020年) in rows 2-12.
# Append the new 2021年 data as row 13.
$ws.Cells.Item(13, 1).Value = "2021年"
$ws.Cells.Item(13, 2).Value = 9172268
$ws.Cells.Item(13, 3).Value = 491036
$ws.Cells.Item(13, 4).Value = 16022746
$ws.Cells.Item(13, 6).Value = 1881576
$ws.Cells.Item(13, 7).Value = 26326333
$ws.Cells.Item(13, 8).Value = 376709
$ws.Cells.Item(13, 9).Value = 1460781
$ws.Cells.Item(13, 10).Value = 5507313
$ws.Cells.Item(13, 11).Value = 278514971
$ws.Cells.Item(13, 12).Value = 9291631
$ws.Cells.Item(13, 13).Value = 36958161
$ws.Cells.Item(13, 14).Value = 273111
$ws.Cells.Item(13, 15).Value = 1201556
$ws.Cells.Item(13, 16).Value = 285419
$ws.Cells.Item(13, 17).Value = 5049240
$ws.Cells.Item(13, 18).Value = 4507518
$ws.Cells.Item(13, 19).Value = 111523784
$ws.Cells.Item(13, 20).Value = 18150765
$ws.Cells.Item(13, 21).Value = 30035025

# Column A on every data row carries the bold/centered/bordered style (style index 1);
# copy that formatting from the row above onto the new A13 cell.
$ws.Range("A12").Copy()
$ws.Range("A13").PasteSpecial(-4122)  # xlPasteFormats

$excel.CutCopyMode = 0

Write-Output "edit applied"
